$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New COUNT (column C) values for the rows whose counts changed.
# Row -> new C value (only rows that actually changed are listed here).
$newCounts = @{
    2  = 522
    5  = 4677
    7  = 1160
    8  = 1236
    14 = 11158
    17 = 75946
    20 = 250
    21 = 1291
    24 = 19508
    26 = 443
    31 = 17191
}

$firstRow = 2
$lastRow = 31

# Apply the updated COUNT values.
foreach ($r in $newCounts.Keys) {
    $ws.Cells.Item($r, 3).Value2 = $newCounts[$r]
}

# Recompute total of column C across all data rows.
$total = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $total = $total + $ws.Cells.Item($r, 3).Value2
}

# Recompute per-"from" group totals (column A groups), used for column E.
$groupTotals = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $fromKey = $ws.Cells.Item($r, 1).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    if ($groupTotals.ContainsKey($fromKey)) {
        $groupTotals[$fromKey] = $groupTotals[$fromKey] + $c
    } else {
        $groupTotals[$fromKey] = $c
    }
}

# Recompute PROP (column D) = COUNT / overall total,
# and PROB (column E) = COUNT / group total, for every data row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $fromKey = $ws.Cells.Item($r, 1).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 4).Value2 = $c / $total
    $ws.Cells.Item($r, 5).Value2 = $c / $groupTotals[$fromKey]
}
